# Weekly fruit/vegetable update: insert 4 new price rows (week of serial
# 44610) just above the existing row-424 block, pushing the old rows
# 424:447 down to 428:451.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the current row 424 (Excel shifts 424:447 -> 428:451
# and extends the used range automatically).
$ws.Rows("424:427").Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112006
$categoria = "Repollo"
$unidad    = "`$/unidad"
$kgUnid    = 1
$clasif    = "Hortaliza"

function Set-RepolloRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen) {
    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $catId
    $ws.Cells.Item($Row, 7).Value  = $categoria
    $ws.Cells.Item($Row, 8).Value  = $Variedad
    $ws.Cells.Item($Row, 9).Value  = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $kgUnid
    $ws.Cells.Item($Row, 18).Value = $clasif
}

Set-RepolloRow 424 44610 "Copenhague"    "Primera" 970  1200 1300 1249 "Región Metropolitana"
Set-RepolloRow 425 44610 "Crespo record" "Primera" 4300 1000 1100 1050 "Región Metropolitana"
Set-RepolloRow 426 44610 "Crespo record" "Segunda" 1960 900  900  900  "Región Metropolitana"
Set-RepolloRow 427 44610 "Morada(o)"     "Primera" 520  1400 1500 1450 "Región Metropolitana"
